$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths
# NOTE: the engine stores ColumnWidth internally with an extra +5/6 (0.8333...)
# padding baked into the XML "width" attribute whenever ColumnWidth is assigned
# via COM, so we subtract that padding here to land exactly on the target
# integer widths from the diff.
$pad = 5 / 6
$ws.Columns.Item(1).ColumnWidth = 66 - $pad
$ws.Columns.Item(2).ColumnWidth = 44 - $pad
$ws.Columns.Item(3).ColumnWidth = 9 - $pad
$ws.Columns.Item(4).ColumnWidth = 45 - $pad
$ws.Columns.Item(5).ColumnWidth = 47 - $pad
$ws.Columns.Item(6).ColumnWidth = 35 - $pad
$ws.Columns.Item(7).ColumnWidth = 37 - $pad
$ws.Columns.Item(8).ColumnWidth = 32 - $pad
$ws.Columns.Item(9).ColumnWidth = 34 - $pad
$ws.Columns.Item(10).ColumnWidth = 40 - $pad
$ws.Columns.Item(11).ColumnWidth = 42 - $pad
$ws.Columns.Item(12).ColumnWidth = 35 - $pad
$ws.Columns.Item(13).ColumnWidth = 37 - $pad

# Update header row text
$ws.Range("B1").Value = "div_testRunComponents_internalRoleCellName"
$ws.Range("D1").Value = "link_projectNavigation_internalRoleLinkName"
$ws.Range("E1").Value = "link_projectNavigation_internalRoleLinkName_1"
$ws.Range("F1").Value = "link_projectNavigation_project_id"
$ws.Range("G1").Value = "link_projectNavigation_project_id_1"
$ws.Range("H1").Value = "link_projectNavigation_team_id"
$ws.Range("I1").Value = "link_projectNavigation_team_id_1"
$ws.Range("J1").Value = "link_projectNavigation_test_project_id"
$ws.Range("K1").Value = "link_projectNavigation_test_project_id_1"
$ws.Range("L1").Value = "link_projectNavigation_trNthChild"
$ws.Range("M1").Value = "link_projectNavigation_trNthChild_1"

# Update data row
$ws.Range("A2").Value = "Data Files/AI-Generated/Common/scheduleAndRunTestSuite-test-data"
